$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 7.863513333333333
$ws.Cells.Item(2, 8).Value = 23.59054
$ws.Cells.Item(2, 9).Value = 0.1316713470554304
$ws.Cells.Item(2, 10).Value = 0.1376659241260802
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 2.476839
$ws.Cells.Item(2, 14).Value = 7.430517
$ws.Cells.Item(2, 15).Value = 0.0525987134655237
$ws.Cells.Item(2, 16).Value = 0.05675564862155354
$ws.Cells.Item(2, 17).Value = 19.47665650102
$ws.Cells.Item(2, 18).Value = 175.28990850918
$ws.Cells.Item(2, 19).Value = 0.006925743455388113
$ws.Cells.Item(2, 20).Value = 0.007813318816861257

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 7.863513333333333
$ws.Cells.Item(3, 8).Value = 23.59054
$ws.Cells.Item(3, 9).Value = 0.1316713470554304
$ws.Cells.Item(3, 10).Value = 0.1376659241260802
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 7.339638666666666
$ws.Cells.Item(3, 14).Value = 22.018916
$ws.Cells.Item(3, 15).Value = 0.1558662275458673
$ws.Cells.Item(3, 16).Value = 0.1681845098427879
$ws.Cells.Item(3, 17).Value = 57.71534651718221
$ws.Cells.Item(3, 18).Value = 519.4381186546399
$ws.Cells.Item(3, 19).Value = 0.02052311614141258
$ws.Cells.Item(3, 20).Value = 0.02315327597119923

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 7.863513333333333
$ws.Cells.Item(4, 8).Value = 23.59054
$ws.Cells.Item(4, 9).Value = 0.1316713470554304
$ws.Cells.Item(4, 10).Value = 0.1376659241260802
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 9.137454333333332
$ws.Cells.Item(4, 14).Value = 27.412363
$ws.Cells.Item(4, 15).Value = 0.1940450478546679
$ws.Cells.Item(4, 16).Value = 0.2093806450230146
$ws.Cells.Item(4, 17).Value = 71.8524939828911
$ws.Cells.Item(4, 18).Value = 646.6724458460199
$ws.Cells.Item(4, 19).Value = 0.02555017284045959
$ws.Cells.Item(4, 20).Value = 0.02882457999120806

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 7.863513333333333
$ws.Cells.Item(5, 8).Value = 23.59054
$ws.Cells.Item(5, 9).Value = 0.1316713470554304
$ws.Cells.Item(5, 10).Value = 0.1376659241260802
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 17.78856566666667
$ws.Cells.Item(5, 14).Value = 53.365697
$ws.Cells.Item(5, 15).Value = 0.3777620057111716
$ws.Cells.Item(5, 16).Value = 0.4076169595435007
$ws.Cells.Item(5, 17).Value = 139.8806233007089
$ws.Cells.Item(5, 18).Value = 1258.92560970638
$ws.Cells.Item(5, 19).Value = 0.04974043215835117
$ws.Cells.Item(5, 20).Value = 0.05611496542501906

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 7.863513333333333
$ws.Cells.Item(6, 8).Value = 23.59054
$ws.Cells.Item(6, 9).Value = 0.1316713470554304
$ws.Cells.Item(6, 10).Value = 0.1376659241260802
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 13).Value = 10.346848
$ws.Cells.Item(6, 14).Value = 20.693696
$ws.Cells.Item(6, 15).Value = 0.2197280054227695
$ws.Cells.Item(6, 16).Value = 0.1580622369691433
$ws.Cells.Item(6, 17).Value = 81.36257720597332
$ws.Cells.Item(6, 18).Value = 488.1754632358399
$ws.Cells.Item(6, 19).Value = 0.02893188245981898
$ws.Cells.Item(6, 20).Value = 0.02175978392179259

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 15.12977766666667
$ws.Cells.Item(7, 8).Value = 45.389333
$ws.Cells.Item(7, 9).Value = 0.2533420014148681
$ws.Cells.Item(7, 10).Value = 0.2648758558689792
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 2.476839
$ws.Cells.Item(7, 14).Value = 7.430517
$ws.Cells.Item(7, 15).Value = 0.0525987134655237
$ws.Cells.Item(7, 16).Value = 0.05675564862155354
$ws.Cells.Item(7, 17).Value = 37.47402338612901
$ws.Cells.Item(7, 18).Value = 337.266210475161
$ws.Cells.Item(7, 19).Value = 0.01332546334120295
$ws.Cells.Item(7, 20).Value = 0.01503320100403304

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 15.12977766666667
$ws.Cells.Item(8, 8).Value = 45.389333
$ws.Cells.Item(8, 9).Value = 0.2533420014148681
$ws.Cells.Item(8, 10).Value = 0.2648758558689792
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 7.339638666666666
$ws.Cells.Item(8, 14).Value = 22.018916
$ws.Cells.Item(8, 15).Value = 0.1558662275458673
$ws.Cells.Item(8, 16).Value = 0.1681845098427879
$ws.Cells.Item(8, 17).Value = 111.0471011803364
$ws.Cells.Item(8, 18).Value = 999.4239106230278
$ws.Cells.Item(8, 19).Value = 0.03948746203945525
$ws.Cells.Item(8, 20).Value = 0.04454801598851321

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 15.12977766666667
$ws.Cells.Item(9, 8).Value = 45.389333
$ws.Cells.Item(9, 9).Value = 0.2533420014148681
$ws.Cells.Item(9, 10).Value = 0.2648758558689792
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 9.137454333333332
$ws.Cells.Item(9, 14).Value = 27.412363
$ws.Cells.Item(9, 15).Value = 0.1940450478546679
$ws.Cells.Item(9, 16).Value = 0.2093806450230146
$ws.Cells.Item(9, 17).Value = 138.2476525026532
$ws.Cells.Item(9, 18).Value = 1244.228872523879
$ws.Cells.Item(9, 19).Value = 0.04915976078814542
$ws.Cells.Item(9, 20).Value = 0.05545987755286991

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 15.12977766666667
$ws.Cells.Item(10, 8).Value = 45.389333
$ws.Cells.Item(10, 9).Value = 0.2533420014148681
$ws.Cells.Item(10, 10).Value = 0.2648758558689792
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 17.78856566666667
$ws.Cells.Item(10, 14).Value = 53.365697
$ws.Cells.Item(10, 15).Value = 0.3777620057111716
$ws.Cells.Item(10, 16).Value = 0.4076169595435007
$ws.Cells.Item(10, 17).Value = 269.1370435455668
$ws.Cells.Item(10, 18).Value = 2422.233391910101
$ws.Cells.Item(10, 19).Value = 0.09570298258536303
$ws.Cells.Item(10, 20).Value = 0.1079678910257958

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 15.12977766666667
$ws.Cells.Item(11, 8).Value = 45.389333
$ws.Cells.Item(11, 9).Value = 0.2533420014148681
$ws.Cells.Item(11, 10).Value = 0.2648758558689792
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 13).Value = 10.346848
$ws.Cells.Item(11, 14).Value = 20.693696
$ws.Cells.Item(11, 15).Value = 0.2197280054227695
$ws.Cells.Item(11, 16).Value = 0.1580622369691433
$ws.Cells.Item(11, 17).Value = 156.5455097907947
$ws.Cells.Item(11, 18).Value = 939.2730587447679
$ws.Cells.Item(11, 19).Value = 0.0556663326607014
$ws.Cells.Item(11, 20).Value = 0.04186687029776724

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 14.38236133333333
$ws.Cells.Item(12, 8).Value = 43.147084
$ws.Cells.Item(12, 9).Value = 0.2408268175206591
$ws.Cells.Item(12, 10).Value = 0.2517908955117437
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 2.476839
$ws.Cells.Item(12, 14).Value = 7.430517
$ws.Cells.Item(12, 15).Value = 0.0525987134655237
$ws.Cells.Item(12, 16).Value = 0.05675564862155354
$ws.Cells.Item(12, 17).Value = 35.622793462492
$ws.Cells.Item(12, 18).Value = 320.605141162428
$ws.Cells.Item(12, 19).Value = 0.01266718076958311
$ws.Cells.Item(12, 20).Value = 0.01429055559177082

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 14.38236133333333
$ws.Cells.Item(13, 8).Value = 43.147084
$ws.Cells.Item(13, 9).Value = 0.2408268175206591
$ws.Cells.Item(13, 10).Value = 0.2517908955117437
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 7.339638666666666
$ws.Cells.Item(13, 14).Value = 22.018916
$ws.Cells.Item(13, 15).Value = 0.1558662275458673
$ws.Cells.Item(13, 16).Value = 0.1681845098427879
$ws.Cells.Item(13, 17).Value = 105.5613353601049
$ws.Cells.Item(13, 18).Value = 950.0520182409439
$ws.Cells.Item(13, 19).Value = 0.03753676753882211
$ws.Cells.Item(13, 20).Value = 0.04234732834451924

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 14.38236133333333
$ws.Cells.Item(14, 8).Value = 43.147084
$ws.Cells.Item(14, 9).Value = 0.2408268175206591
$ws.Cells.Item(14, 10).Value = 0.2517908955117437
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 9.137454333333332
$ws.Cells.Item(14, 14).Value = 27.412363
$ws.Cells.Item(14, 15).Value = 0.1940450478546679
$ws.Cells.Item(14, 16).Value = 0.2093806450230146
$ws.Cells.Item(14, 17).Value = 131.4181698888324
$ws.Cells.Item(14, 18).Value = 1182.763528999492
$ws.Cells.Item(14, 19).Value = 0.04673125133048368
$ws.Cells.Item(14, 20).Value = 0.05272014011317135

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 14.38236133333333
$ws.Cells.Item(15, 8).Value = 43.147084
$ws.Cells.Item(15, 9).Value = 0.2408268175206591
$ws.Cells.Item(15, 10).Value = 0.2517908955117437
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 17.78856566666667
$ws.Cells.Item(15, 14).Value = 53.365697
$ws.Cells.Item(15, 15).Value = 0.3777620057111716
$ws.Cells.Item(15, 16).Value = 0.4076169595435007
$ws.Cells.Item(15, 17).Value = 255.8415790197276
$ws.Cells.Item(15, 18).Value = 2302.574211177548
$ws.Cells.Item(15, 19).Value = 0.0909752216156425
$ws.Cells.Item(15, 20).Value = 0.1026342392692322

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 14.38236133333333
$ws.Cells.Item(16, 8).Value = 43.147084
$ws.Cells.Item(16, 9).Value = 0.2408268175206591
$ws.Cells.Item(16, 10).Value = 0.2517908955117437
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 13).Value = 10.346848
$ws.Cells.Item(16, 14).Value = 20.693696
$ws.Cells.Item(16, 15).Value = 0.2197280054227695
$ws.Cells.Item(16, 16).Value = 0.1580622369691433
$ws.Cells.Item(16, 17).Value = 148.8121065970773
$ws.Cells.Item(16, 18).Value = 892.872639582464
$ws.Cells.Item(16, 19).Value = 0.0529163962661277
$ws.Cells.Item(16, 20).Value = 0.03979863219305003

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 14.54360033333333
$ws.Cells.Item(17, 8).Value = 43.630801
$ws.Cells.Item(17, 9).Value = 0.243526699294608
$ws.Cells.Item(17, 10).Value = 0.2546136943039924
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 2.476839
$ws.Cells.Item(17, 14).Value = 7.430517
$ws.Cells.Item(17, 15).Value = 0.0525987134655237
$ws.Cells.Item(17, 16).Value = 0.05675564862155354
$ws.Cells.Item(17, 17).Value = 36.022156506013
$ws.Cells.Item(17, 18).Value = 324.199408554117
$ws.Cells.Item(17, 19).Value = 0.01280919107740184
$ws.Cells.Item(17, 20).Value = 0.01445076536815304

# Row 18
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 7).Value = 14.54360033333333
$ws.Cells.Item(18, 8).Value = 43.630801
$ws.Cells.Item(18, 9).Value = 0.243526699294608
$ws.Cells.Item(18, 10).Value = 0.2546136943039924
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 13).Value = 7.339638666666666
$ws.Cells.Item(18, 14).Value = 22.018916
$ws.Cells.Item(18, 15).Value = 0.1558662275458673
$ws.Cells.Item(18, 16).Value = 0.1681845098427879
$ws.Cells.Item(18, 17).Value = 106.7447713590795
$ws.Cells.Item(18, 18).Value = 960.7029422317158
$ws.Cells.Item(18, 19).Value = 0.03795758792574735
$ws.Cells.Item(18, 20).Value = 0.04282207937577841

# Row 19
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 7).Value = 14.54360033333333
$ws.Cells.Item(19, 8).Value = 43.630801
$ws.Cells.Item(19, 9).Value = 0.243526699294608
$ws.Cells.Item(19, 10).Value = 0.2546136943039924
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 13).Value = 9.137454333333332
$ws.Cells.Item(19, 14).Value = 27.412363
$ws.Cells.Item(19, 15).Value = 0.1940450478546679
$ws.Cells.Item(19, 16).Value = 0.2093806450230146
$ws.Cells.Item(19, 17).Value = 132.8914838880848
$ws.Cells.Item(19, 18).Value = 1196.023354992763
$ws.Cells.Item(19, 19).Value = 0.04725515001851153
$ws.Cells.Item(19, 20).Value = 0.05331117954506258

# Row 20
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 7).Value = 14.54360033333333
$ws.Cells.Item(20, 8).Value = 43.630801
$ws.Cells.Item(20, 9).Value = 0.243526699294608
$ws.Cells.Item(20, 10).Value = 0.2546136943039924
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 13).Value = 17.78856566666667
$ws.Cells.Item(20, 14).Value = 53.365697
$ws.Cells.Item(20, 15).Value = 0.3777620057111716
$ws.Cells.Item(20, 16).Value = 0.4076169595435007
$ws.Cells.Item(20, 17).Value = 258.7097895592552
$ws.Cells.Item(20, 18).Value = 2328.388106033297
$ws.Cells.Item(20, 19).Value = 0.09199513436975246
$ws.Cells.Item(20, 20).Value = 0.1037848599303317

# Row 21
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 7).Value = 14.54360033333333
$ws.Cells.Item(21, 8).Value = 43.630801
$ws.Cells.Item(21, 9).Value = 0.243526699294608
$ws.Cells.Item(21, 10).Value = 0.2546136943039924
$ws.Cells.Item(21, 11).Value = 2
$ws.Cells.Item(21, 13).Value = 10.346848
$ws.Cells.Item(21, 14).Value = 20.693696
$ws.Cells.Item(21, 15).Value = 0.2197280054227695
$ws.Cells.Item(21, 16).Value = 0.1580622369691433
$ws.Cells.Item(21, 17).Value = 150.4804220217493
$ws.Cells.Item(21, 18).Value = 902.882532130496
$ws.Cells.Item(21, 19).Value = 0.05350963590319477
$ws.Cells.Item(21, 20).Value = 0.04024481008466666

# Row 22
$ws.Cells.Item(22, 5).Value = 2
$ws.Cells.Item(22, 7).Value = 7.8015105
$ws.Cells.Item(22, 8).Value = 15.603021
$ws.Cells.Item(22, 9).Value = 0.1306331347144344
$ws.Cells.Item(22, 10).Value = 0.0910536301892045
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 13).Value = 2.476839
$ws.Cells.Item(22, 14).Value = 7.430517
$ws.Cells.Item(22, 15).Value = 0.0525987134655237
$ws.Cells.Item(22, 16).Value = 0.05675564862155354
$ws.Cells.Item(22, 17).Value = 19.3230854653095
$ws.Cells.Item(22, 18).Value = 115.938512791857
$ws.Cells.Item(22, 19).Value = 0.006871134821947694
$ws.Cells.Item(22, 20).Value = 0.005167807840735369

# Row 23
$ws.Cells.Item(23, 5).Value = 2
$ws.Cells.Item(23, 7).Value = 7.8015105
$ws.Cells.Item(23, 8).Value = 15.603021
$ws.Cells.Item(23, 9).Value = 0.1306331347144344
$ws.Cells.Item(23, 10).Value = 0.0910536301892045
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 13).Value = 7.339638666666666
$ws.Cells.Item(23, 14).Value = 22.018916
$ws.Cells.Item(23, 15).Value = 0.1558662275458673
$ws.Cells.Item(23, 16).Value = 0.1681845098427879
$ws.Cells.Item(23, 17).Value = 57.26026812420599
$ws.Cells.Item(23, 18).Value = 343.561608745236
$ws.Cells.Item(23, 19).Value = 0.02036129390042997
$ws.Cells.Item(23, 20).Value = 0.01531381016277784

# Row 24
$ws.Cells.Item(24, 5).Value = 2
$ws.Cells.Item(24, 7).Value = 7.8015105
$ws.Cells.Item(24, 8).Value = 15.603021
$ws.Cells.Item(24, 9).Value = 0.1306331347144344
$ws.Cells.Item(24, 10).Value = 0.0910536301892045
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 13).Value = 9.137454333333332
$ws.Cells.Item(24, 14).Value = 27.412363
$ws.Cells.Item(24, 15).Value = 0.1940450478546679
$ws.Cells.Item(24, 16).Value = 0.2093806450230146
$ws.Cells.Item(24, 17).Value = 71.28594592477049
$ws.Cells.Item(24, 18).Value = 427.715675548623
$ws.Cells.Item(24, 19).Value = 0.02534871287706771
$ws.Cells.Item(24, 20).Value = 0.01906486782070267

# Row 25
$ws.Cells.Item(25, 5).Value = 2
$ws.Cells.Item(25, 7).Value = 7.8015105
$ws.Cells.Item(25, 8).Value = 15.603021
$ws.Cells.Item(25, 9).Value = 0.1306331347144344
$ws.Cells.Item(25, 10).Value = 0.0910536301892045
$ws.Cells.Item(25, 11).Value = 3
$ws.Cells.Item(25, 13).Value = 17.78856566666667
$ws.Cells.Item(25, 14).Value = 53.365697
$ws.Cells.Item(25, 15).Value = 0.3777620057111716
$ws.Cells.Item(25, 16).Value = 0.4076169595435007
$ws.Cells.Item(25, 17).Value = 138.7776818284395
$ws.Cells.Item(25, 18).Value = 832.666090970637
$ws.Cells.Item(25, 19).Value = 0.04934823498206243
$ws.Cells.Item(25, 20).Value = 0.03711500389312184

# Row 26
$ws.Cells.Item(26, 5).Value = 2
$ws.Cells.Item(26, 7).Value = 7.8015105
$ws.Cells.Item(26, 8).Value = 15.603021
$ws.Cells.Item(26, 9).Value = 0.1306331347144344
$ws.Cells.Item(26, 10).Value = 0.0910536301892045
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 13).Value = 10.346848
$ws.Cells.Item(26, 14).Value = 20.693696
$ws.Cells.Item(26, 15).Value = 0.2197280054227695
$ws.Cells.Item(26, 16).Value = 0.1580622369691433
$ws.Cells.Item(26, 17).Value = 80.72104331390399
$ws.Cells.Item(26, 18).Value = 322.884173255616
$ws.Cells.Item(26, 19).Value = 0.02870375813292662
$ws.Cells.Item(26, 20).Value = 0.01439214047186678
